# TC24_Canine_Filter_Breed-FrenchBullDg.xlsx
#
# Commit: "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The "CasesTab" Cypher query (cell B2 on the "startup" sheet) incorrectly
# returned an extra `Cohort` column. Fix it by removing the trailing
# `coalesce(co.cohort_description, '') AS `Cohort`` line (and the now
# trailing comma on the previous line) from the RETURN clause.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$fixedQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
  "WHERE demo.breed IN ['French Bulldog']`n" +
  "MATCH (c)<--(diag:diagnosis)`n" +
  "OPTIONAL MATCH (samp:sample)-->(c)`n" +
  "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
  "WITH DISTINCT c, s, demo, diag, co`n" +
  "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
  "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
  "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
  "        coalesce(demo.breed, '') AS Breed ,`n" +
  "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
  "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
  "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
  "        coalesce(demo.sex, '') AS Sex ,`n" +
  "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
  "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
  "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $fixedQuery

# Row 2 shrank (one fewer RETURN line) while rows 3/4 text is unchanged but
# still re-wraps a touch differently - match the author's resulting row
# heights.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 216
$ws.Rows.Item(4).RowHeight = 244.8

# The author's selection ended up on B2 (the cell they just edited) rather
# than the previous C4/topLeftCell A4 scroll position.
$ws.Activate()
$ws.Range("B2").Select()
